$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting existing rows 84..183 down to 85..184
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly data point
$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 44966
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112030
$ws.Range("G84").Value = "Poroto granado"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 400
$ws.Range("K84").Value = 25000
$ws.Range("L84").Value = 27000
$ws.Range("M84").Value = 26000
$ws.Range("N84").Value = "$/saco 25 kilos"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 1040
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
